# Add a new 4th row of reaction-time data to the sheet, mirroring the
# existing rows' layout (Pessoa / Alvo / idioma / lista / Tempo Reaçao em segundos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "9135C37F"
$ws.Range("B4").Value = "汉语"
$ws.Range("C4").Value = "Japones"
$ws.Range("D4").Value = "A"
# Leading apostrophe forces Excel to store this numeric-looking value ("0")
# as text instead of silently coercing it to a number.
$ws.Range("E4").Value = "'0"
